# Apply the "cryptos list" refresh (prices / 1h volume%, plus the
# Toncoin<->Cardano and Bittensor<->PancakeSwap row swaps).
#
# Column D holds price strings that *look* numeric ("29.40", "1.00", ...).
# A plain Range.Value assignment lets Excel's COM layer auto-coerce those
# into real numbers (dropping the original text formatting, e.g. trailing
# zeros). To keep them as literal text we briefly force the cell to the
# "@" (Text) number format, assign the string, then reset the cell's
# style back to "Normal" so no stray style index is left attached to the
# cell (matches the original workbook, where none of these cells carry
# an explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '65.735.58'
$cell.Style = 'Normal'

$ws.Range('E2').Value = '  -0.03%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '2.679.33'
$cell.Style = 'Normal'

$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('E4').Value = '  -0.04%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '600.75'
$cell.Style = 'Normal'

$ws.Range('E5').Value = '  -0.89%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '156.93'
$cell.Style = 'Normal'

$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  -0.06%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.622'
$cell.Style = 'Normal'

$ws.Range('E8').Value = '  +6.02%  '
$ws.Range('E9').Value = '  +5.23%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.401'
$cell.Style = 'Normal'

$ws.Range('E10').Value = '  -0.15%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '5.89'
$cell.Style = 'Normal'

$ws.Range('E11').Value = '  -2.69%  '
$ws.Range('E12').Value = '  -0.05%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '29.40'
$cell.Style = 'Normal'

$ws.Range('E13').Value = '  -2.21%  '
$ws.Range('E14').Value = '  -2.14%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '3.159.68'
$cell.Style = 'Normal'

$ws.Range('E15').Value = '  -0.85%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '65.612.07'
$cell.Style = 'Normal'

$ws.Range('E16').Value = '  -0.06%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '2.669.79'
$cell.Style = 'Normal'

$ws.Range('E17').Value = '  -1.79%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '12.91'
$cell.Style = 'Normal'

$ws.Range('E18').Value = '  +1.56%  '
$ws.Range('E19').Value = '  -1.26%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '7.61'
$cell.Style = 'Normal'

$ws.Range('E20').Value = '  +1.42%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '352.41'
$cell.Style = 'Normal'

$ws.Range('E21').Value = '  -2.06%  '
$ws.Range('E22').Value = '  +0.00%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '69.80'
$cell.Style = 'Normal'

$ws.Range('E23').Value = '  -0.79%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '0.0000112'
$cell.Style = 'Normal'

$ws.Range('E24').Value = '  +6.02%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '9.68'
$cell.Style = 'Normal'

$ws.Range('E25').Value = '  -1.58%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  -2.02%  '
$ws.Range('E28').Value = '  -5.56%  '
$ws.Range('E29').Value = '  -2.82%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'

$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '2.14'
$cell.Style = 'Normal'

$ws.Range('E31').Value = '  -2.59%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '529.33'
$cell.Style = 'Normal'

$ws.Range('E32').Value = '  -2.01%  '
$ws.Range('E33').Value = '  -1.40%  '
$ws.Range('E34').Value = '  -2.89%  '
$ws.Range('E35').Value = '  +1.95%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.425'
$cell.Style = 'Normal'

$ws.Range('E36').Value = '  -1.73%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '20.70'
$cell.Style = 'Normal'

$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('E38').Value = '  -0.01%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '158.08'
$cell.Style = 'Normal'

$ws.Range('E39').Value = '  -2.77%  '
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('E41').Value = '  +0.00%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '164.42'
$cell.Style = 'Normal'

$ws.Range('E43').Value = '  -1.02%  '
$ws.Range('E44').Value = '  +2.92%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.0612'
$cell.Style = 'Normal'

$ws.Range('E45').Value = '  -0.71%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '22.94'
$cell.Style = 'Normal'

$ws.Range('E46').Value = '  -2.80%  '
$ws.Range('E47').Value = '  +17.55%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '0.0260'
$cell.Style = 'Normal'

$ws.Range('E48').Value = '  -2.43%  '
$ws.Range('E49').Value = '  -2.56%  '
$ws.Range('E50').Value = '  +2.42%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '20.18'
$cell.Style = 'Normal'

$ws.Range('E51').Value = '  -4.21%  '
